# Updated cryptos list data values (prices / 1h volume %) to match upstream scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.496.77'
$ws.Range("E2").Value = '  -1.13%  '
$ws.Range("D3").Value = '3.761.49'
$ws.Range("E3").Value = '  -1.51%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.60'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.64'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.77%  '
$ws.Range("D7").Value = '3.761.47'
$ws.Range("E7").Value = '  -1.42%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  -0.12%  '
$ws.Range("E10").Value = '  +1.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.54'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.93%  '
$ws.Range("E12").Value = '  -0.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000277'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +6.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.68'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.87%  '
$ws.Range("D15").Value = '4.393.54'
$ws.Range("E15").Value = '  -1.31%  '
$ws.Range("D16").Value = '3.762.36'
$ws.Range("E16").Value = '  -1.17%  '
$ws.Range("E17").Value = '  +0.47%  '
$ws.Range("D18").Value = '67.520.99'
$ws.Range("E18").Value = '  -0.95%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.22'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.53%  '
$ws.Range("E20").Value = '  +1.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.56'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.75%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '469.04'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.721'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.70%  '
$ws.Range("E24").Value = '  -8.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.92'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.33%  '
$ws.Range("E26").Value = '  -0.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.19'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.49%  '
$ws.Range("E28").Value = '  +3.67%  '
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("E30").Value = '  -1.61%  '
$ws.Range("D31").Value = '3.909.77'
$ws.Range("E31").Value = '  -1.29%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.65'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '30.65'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.53%  '
$ws.Range("E34").Value = '  -2.86%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.16'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -3.29%  '
$ws.Range("D36").Value = '3.727.58'
$ws.Range("E36").Value = '  -1.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.87'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +5.02%  '
$ws.Range("E38").Value = '  +0.31%  '
$ws.Range("E39").Value = '  -1.78%  '
$ws.Range("E40").Value = '  -1.47%  '
$ws.Range("E41").Value = '  -1.08%  '
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.313'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.59%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.75'
$ws.Range("D45").ClearFormats()
$ws.Range("E46").Value = '  -1.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '45.91'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '399.83'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -4.43%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000270'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -8.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '140.20'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.89%  '
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0355'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.86%  '
